$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-03-17 Sunday" "2024-03-18 Monday"
Replace-Text "770÷4=" "227÷6="
Replace-Text "591÷6=" "923÷6="
Replace-Text "524÷4=" "731÷4="
Replace-Text "864÷2=" "145÷7="
Replace-Text "834÷7=" "369÷2="
Replace-Text "305÷2=" "373÷3="
Replace-Text "371÷9=" "475÷5="
Replace-Text "975÷4=" "622÷3="
Replace-Text "735÷2=" "395÷3="
Replace-Text "864÷4=" "509÷2="
Replace-Text "674÷7=" "864÷7="
Replace-Text "480÷8=" "200÷8="
Replace-Text "830÷3=" "658÷9="
Replace-Text "723÷2=" "959÷8="
Replace-Text "165÷9=" "127÷4="
Replace-Text "225÷3=" "896÷3="
Replace-Text "693÷9=" "758÷9="
Replace-Text "144÷8=" "999÷4="
Replace-Text "209÷7=" "167÷3="
Replace-Text "232÷4=" "295÷9="
Replace-Text "599÷6=" "372÷7="
Replace-Text "759÷2=" "950÷2="
Replace-Text "842÷9=" "927÷9="
Replace-Text "911÷7=" "732÷9="
Replace-Text "206÷6=" "123÷9="
